$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9413926005363464
$ws.Range("B1").Value = 1.976456165313721
$ws.Range("C1").Value = 7.435521125793457
$ws.Range("D1").Value = 2.725505352020264
$ws.Range("E1").Value = 1.423065304756165
